$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at the top of the data block (rows 2-7), pushing the
# existing accelerometer readings down by 6 rows.
$ws.Rows("2:7").Insert()
# The inserted rows pick up the bold header formatting from row 1 above;
# clear that so the new data cells stay unstyled, like the rest of the data.
$ws.Range("A2:C7").ClearFormats()

# New readings recorded at the start of the series.
$newTop = @(
    @(-2.778316736221313, 6.619067668914795, -0.6493567824363708),
    @(-2.981808757781982, 6.65326156616211, -0.8340041637420658),
    @(-2.862793350219726, 6.595655870437622, -1.086881220340729),
    @(-2.658387470245361, 6.523494625091553, -1.100642728805541),
    @(-2.623731708526611, 6.502191925048828, -0.9057361066341401),
    @(-2.912428379058838, 6.643356800079346, -0.9624049067497252)
)

$r = 2
foreach ($row in $newTop) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# New readings appended at the end of the series (rows 28-31).
$newBottom = @(
    @(1.882777214050291, 5.883625364303589, -1.771220207214356),
    @(2.030305290222169, 6.322917270660403, -1.870007395744325),
    @(2.015394306182861, 6.204385328292845, -1.822170174121856),
    @(1.990867900848389, 6.087325572967529, -1.848094010353089)
)

$r = 28
foreach ($row in $newBottom) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
